$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text-typed so numeric-looking price strings (e.g. "316.20")
# are not auto-converted to numbers, matching the original inlineStr cell type.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "24.532.09"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "1.698.11"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").Value = "316.20"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").Value = "0.3904"
$ws.Range("E7").Value = "  -0.66%  "
$ws.Range("D8").Value = "0.4084"
$ws.Range("E8").Value = "  +1.22%  "
$ws.Range("E9").Value = "  -1.96%  "
$ws.Range("D10").Value = "1.001"
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("D11").Value = "52.82"
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("D12").Value = "0.08813"
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").Value = "26.50"
$ws.Range("E13").Value = "  +12.21%  "
$ws.Range("D14").Value = "7.511"
$ws.Range("E14").Value = "  +1.44%  "
$ws.Range("D15").Value = "8.327"
$ws.Range("E15").Value = "  +2.38%  "
$ws.Range("D16").Value = "0.00001352"
$ws.Range("E16").Value = "  +2.32%  "
$ws.Range("D17").Value = "1.685.39"
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("D18").Value = "98.00"
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("D19").Value = "0.07195"
$ws.Range("E19").Value = "  +2.19%  "
$ws.Range("D20").Value = "20.69"
$ws.Range("E20").Value = "  +5.21%  "
$ws.Range("D21").Value = "7.331"
$ws.Range("E21").Value = "  +3.58%  "
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  -0.64%  "
$ws.Range("D23").Value = "14.41"
$ws.Range("E23").Value = "  -2.39%  "
$ws.Range("D24").Value = "24.537.98"
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("D25").Value = "3.045"
$ws.Range("E25").Value = "  -3.03%  "
$ws.Range("E26").Value = "  -1.44%  "
$ws.Range("D27").Value = "23.15"
$ws.Range("E27").Value = "  +2.11%  "
$ws.Range("D28").Value = "168.35"
$ws.Range("E28").Value = "  +3.20%  "
$ws.Range("D29").Value = "146.74"
$ws.Range("E29").Value = "  +8.16%  "
$ws.Range("D30").Value = "8.485"
$ws.Range("E30").Value = "  -2.92%  "
$ws.Range("D31").Value = "5.399"
$ws.Range("E31").Value = "  +4.20%  "
$ws.Range("D32").Value = "1.876.26"
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("D33").Value = "2.205"
$ws.Range("E33").Value = "  +11.56%  "
$ws.Range("D34").Value = "0.08825"
$ws.Range("E34").Value = "  -2.45%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "1.057"
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "7.279"
$ws.Range("E36").Value = "  -5.23%  "
$ws.Range("D37").Value = "0.03116"
$ws.Range("E37").Value = "  +12.69%  "
$ws.Range("D38").Value = "0.2812"
$ws.Range("E38").Value = "  +2.14%  "
$ws.Range("D39").Value = "10.99"
$ws.Range("E39").Value = "  -0.74%  "
$ws.Range("D40").Value = "0.09196"
$ws.Range("E40").Value = "  +0.79%  "
$ws.Range("D41").Value = "14.29"
$ws.Range("E41").Value = "  -1.64%  "
$ws.Range("D42").Value = "0.8015"
$ws.Range("E42").Value = "  +4.49%  "
$ws.Range("D43").Value = "1.486"
$ws.Range("E43").Value = "  +1.61%  "
$ws.Range("D44").Value = "17.49"
$ws.Range("E44").Value = "  +9.79%  "
$ws.Range("D45").Value = "2.687"
$ws.Range("E45").Value = "  +4.90%  "
$ws.Range("D46").Value = "0.7292"
$ws.Range("E46").Value = "  +1.83%  "
$ws.Range("D47").Value = "4.270"
$ws.Range("E47").Value = "  +1.46%  "
$ws.Range("E48").Value = "  +6.39%  "
$ws.Range("D49").Value = "1.000"
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("D50").Value = "141.22"
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("D51").Value = "0.08174"
$ws.Range("E51").Value = "  +2.45%  "

# Restore default cell style on column D (NumberFormat change above bumps the
# style index; reset back to Normal so style indices match the original).
$priceRange.Style = "Normal"
